$wb = $excel.ActiveWorkbook

# New business-day dates appended to every sheet (rows 106-110), shared
# across all sheets in the workbook.
$newDates = @(45971, 45972, 45973, 45974, 45975)

# Per-sheet data: the corrected value for the existing placeholder row
# (A105 / "2025-11-07", previously 0) plus the 5 new "remn_amt" values
# that go with $newDates, in sheet order (sheet1 .. sheet7).
$sheetFix = @{
    1 = 3221
    2 = 1146
    3 = 1245
    4 = 1848
    5 = 731
    6 = 1500
    7 = 2845
}
$sheetNewVals = @{
    1 = @(3241, 3212, 3231, 3222, 3254)
    2 = @(1141, 1141, 1159, 1155, 1134)
    3 = @(1394, 1352, 1364, 1336, 1324)
    4 = @(1869, 1872, 1893, 1897, 1902)
    5 = @(750, 741, 759, 762, 762)
    6 = @(1505, 1502, 1508, 1524, 1524)
    7 = @(2627, 2591, 2672, 2685, 2657)
}

for ($i = 1; $i -le 7; $i++) {
    $ws = $wb.Worksheets.Item($i)

    # Fix the previously-placeholder B105 value.
    $ws.Range("B105").Value = $sheetFix[$i]

    # Append rows 106-110: date in column A (matching the existing date
    # formatting), remn_amt in column B.
    $bvals = $sheetNewVals[$i]
    for ($r = 0; $r -lt $newDates.Count; $r++) {
        $rowNum = 106 + $r
        $ws.Range("A$rowNum").Value = $newDates[$r]
        $ws.Range("A$rowNum").NumberFormat = "YYYY-MM-DD HH:MM:SS"
        $ws.Range("B$rowNum").Value = $bvals[$r]
    }
}
